$wb = $excel.ActiveWorkbook

# New validation list (with "seva" appended) used across all repository-source sheets.
# Quoted, matching Excel's native literal-list formula syntax ("a,b,c").
$repoList = '"addgene,genbank,benchling,snapgene,euroscarf,igem,wekwikgene,seva"'

# 1. Insert the new "SEVASource" sheet right after "WekWikGeneIdSource".
$afterSheet = $wb.Worksheets.Item("WekWikGeneIdSource")
$seva = $wb.Worksheets.Add($null, $afterSheet)
$seva.Name = "SEVASource"

# Header row matching the other *IdSource sheets (sequence_file_url .. id).
$headers = @("sequence_file_url", "repository_id", "repository_name", "input", "output", "type", "output_name", "id")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $seva.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# repository_name (column C) gets the dropdown validation, seva included.
$seva.Range("C2:C1048576").Validation.Add(3, 1, 1, $repoList)

# 2. Update the existing repository-name validation lists so they also allow "seva".
$sheetsToUpdate = @(
    @{ Name = "RepositoryIdSource"; Column = "B" },
    @{ Name = "AddGeneIdSource"; Column = "D" },
    @{ Name = "WekWikGeneIdSource"; Column = "C" },
    @{ Name = "BenchlingUrlSource"; Column = "B" },
    @{ Name = "SnapGenePlasmidSource"; Column = "B" },
    @{ Name = "EuroscarfSource"; Column = "B" },
    @{ Name = "IGEMSource"; Column = "C" }
)

foreach ($entry in $sheetsToUpdate) {
    $ws = $wb.Worksheets.Item($entry.Name)
    $rng = $ws.Range("$($entry.Column)2:$($entry.Column)1048576")
    $rng.Validation.Modify(3, 1, 1, $repoList)
}

Write-Output "SEVASource inserted and validations updated"
